$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Meetrapport titel", $true, $false, $false, $false, $false,
    $true, 1, $false, "Welk web server pakket is het beste voor ons?", 2
)
